$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 changes: Z7 1 -> 0, AD7 13 -> 12 ---
$ws.Range("Z7").Value = 0
$ws.Range("AD7").Value = 12

# --- AD column bulk value updates (rows 25-221) ---
$ws.Range("AD25:AD27").Value = 0.1388888888833333
$ws.Range("AD28:AD41").Value = 0.2499999999916667
$ws.Range("AD42:AD73").Value = 0.8333333333250001
$ws.Range("AD74:AD80").Value = 0.645833333325
$ws.Range("AD81:AD87").Value = 0.6291666666583333
$ws.Range("AD88:AD105").Value = 0.4347222222166667
$ws.Range("AD106:AD125").Value = 0.1666666666666667
$ws.Range("AD126:AD221").Value = 0.1944444444416667

# --- New rows 222-233: add 12 additional days of data ---
$newDates = @("9/30/2020", "10/1/2020", "10/2/2020", "10/3/2020", "10/4/2020", "10/5/2020", "10/6/2020", "10/7/2020", "10/8/2020", "10/9/2020", "10/10/2020", "10/11/2020")
$bcVals = @(0, 0, 1, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 1, 1, 0, 1, 0)
$startRow = 222
$endRow = 233

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Formula = '="' + $newDates[$i] + '"'
    $arr = New-Object 'object[,]' 1,28
    for ($j = 0; $j -lt $bcVals.Length; $j++) {
        $arr[0,$j] = $bcVals[$j]
    }
    $ws.Range("B" + $r + ":AC" + $r).Value = $arr
    $ws.Range("AD" + $r).Value = 0.1944444444416667
}

# Convert column-A formulas to static values (avoid date auto-conversion side effects)
$ws.Range("A" + $startRow + ":A" + $endRow).Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4163)

# Copy the style (bold font + border) from A221 onto the new A cells
$ws.Range("A221").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)

$excel.CutCopyMode = 0
